$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 188
$ws.Range("F3").Value = 403
$ws.Range("F4").Value = 1136
$ws.Range("F7").Value = 22
$ws.Range("F8").Value = 1060
$ws.Range("F10").Value = 333
$ws.Range("F11").Value = 417
$ws.Range("F12").Value = 33
$ws.Range("F13").Value = 310
$ws.Range("F14").Value = 351
$ws.Range("F15").Value = 27
$ws.Range("F17").Value = 444
$ws.Range("F18").Value = 441
$ws.Range("F19").Value = 5557
$ws.Range("F21").Value = 1552
$ws.Range("F22").Value = 363
$ws.Range("F23").Value = 4723
$ws.Range("F24").Value = 118
$ws.Range("F25").Value = 84
$ws.Range("F26").Value = 1488
$ws.Range("F29").Value = 646
$ws.Range("F30").Value = 51

# Sheet: 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 8
$ws.Range("F5").Value = 127
$ws.Range("F8").Value = 86
$ws.Range("F18").Value = 1

# Sheet: 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 580
$ws.Range("F4").Value = 2123

# Sheet: 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 580
$ws.Range("F4").Value = 2123
$ws.Range("F5").Value = 188
$ws.Range("F6").Value = 403
$ws.Range("F7").Value = 1136
$ws.Range("F10").Value = 22
$ws.Range("F11").Value = 1060
$ws.Range("F12").Value = 333
$ws.Range("F13").Value = 417
$ws.Range("F14").Value = 33
$ws.Range("F15").Value = 310
$ws.Range("F16").Value = 351
$ws.Range("F17").Value = 27
$ws.Range("F22").Value = 441
$ws.Range("F23").Value = 5557
$ws.Range("F25").Value = 1553
$ws.Range("F28").Value = 363
$ws.Range("F31").Value = 4723
$ws.Range("F32").Value = 118
$ws.Range("F33").Value = 84
$ws.Range("F34").Value = 1488
$ws.Range("F37").Value = 646
$ws.Range("F38").Value = 51
$ws.Range("F45").Value = 1
